$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new numeric value, per the target diff
$changes = @{
    "E3"  = 12.1448
    "A4"  = -21.47170000000002
    "D4"  = -8.149400000000009
    "E4"  = 12.36130000000001
    "D5"  = -8.6934
    "A6"  = -20.09699999999999
    "D6"  = -8.828499999999995
    "A7"  = -21.02100000000001
    "A8"  = -20.73059999999999
    "D8"  = -8.663499999999996
    "E9"  = 13.74820000000001
    "E11" = 13.3062
    "E14" = 13.65540000000001
    "A16" = -20.31559999999999
    "D16" = -8.185800000000002
    "E18" = 12.4041
    "A20" = -22.06620000000003
    "A21" = -20.42989999999999
    "D22" = -8.109699999999997
    "E25" = 13.305
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
